$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 19 blank rows above the old row 11 ("konVAR1" row), pushing every
# row below (old rows 11-50) down to rows 30-69.
$ws.Range("A11:E29").EntireRow.Insert()

# The newly inserted rows come back with default formatting; copy the
# visual style (borders / number format) from the row directly above
# (row 10, which carries the same "body row" style as every other data
# row) onto the freshly inserted block.
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column B: VEGF-A complex concentrations (entered first) ---
$ws.Range("B11").Value = "Concentration of VEGF-A:VEGFR1"
$ws.Range("B12").Value = "Concentration of VEGF-A:VEGFR2"
$ws.Range("B13").Value = "Concentration of VEGF-A:NRP1"
$ws.Range("B14").Value = "Concentration of VEGF-A:VEGFR2:NRP1"
$ws.Range("B15").Value = "Concentration of VEGF-A:PDGFR" + [char]945
$ws.Range("B16").Value = "Concentration of VEGF-A:PDGFR" + [char]946

# --- Column A: short parameter codes for all 19 new rows ---
$ws.Range("A11").Value = "VA_R1"
$ws.Range("A12").Value = "VA_R2"
$ws.Range("A13").Value = "VA_N1"
$ws.Range("A14").Value = "VA_R2_N1"
$ws.Range("A15").Value = "VA_PDRa"
$ws.Range("A16").Value = "VA_PDRb"
$ws.Range("A17").Value = "VB_R1"
$ws.Range("A18").Value = "VB_N1"
$ws.Range("A19").Value = "Pl_R1"
$ws.Range("A20").Value = "Pl_N1"
$ws.Range("A21").Value = "PDAA_R2"
$ws.Range("A22").Value = "PDAA_PDRa"
$ws.Range("A23").Value = "PDAB_R2"
$ws.Range("A24").Value = "PDAB_PDRa"
$ws.Range("A25").Value = "PDAB_PDRb"
$ws.Range("A26").Value = "PDBB_R2"
$ws.Range("A27").Value = "PDBB_PDRa"
$ws.Range("A28").Value = "PDBB_PDRb"
$ws.Range("A29").Value = "R1_N1"

# --- Column B: remaining complex concentrations (entered after column A) ---
$ws.Range("B17").Value = "Concentration of VEGF-B:VEGFR1"
$ws.Range("B18").Value = "Concentration of VEGF-B:PlGF"
$ws.Range("B19").Value = "Concentration of PlGF:VEGFR1"
$ws.Range("B20").Value = "Concentration of PlGF:NRP1"
$ws.Range("B21").Value = "Concentration of PDGF-AA:VEGFR2"
$ws.Range("B22").Value = "Concentration of PDGF-AA:PDGFR" + [char]945
$ws.Range("B23").Value = "Concentration of PDGF-AB:VEGFR2"
$ws.Range("B24").Value = "Concentration of PDGF-AB:PDGFR" + [char]945
$ws.Range("B25").Value = "Concentration of PDGF-AB:PDGFR" + [char]946
$ws.Range("B26").Value = "Concentration of PDGF-BB:VEGFR2"
$ws.Range("B27").Value = "Concentration of PDGF-BB:PDGFR" + [char]945
$ws.Range("B28").Value = "Concentration of PDGF-BB:PDGFR" + [char]946
$ws.Range("B29").Value = "Concentration of VEGFR1:NRP1"

# --- Columns C (value), D (unit) and E (reference) for all 19 new rows ---
$ws.Range("C11:C29").Value = 0
$ws.Range("D11:D29").Value = "M"
$ws.Range("E11:E29").Value = "Assumed"

# Update the view to match the saved window state (scroll position is a
# cosmetic, host-managed property; best effort only).
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
$ws.Range("B15").Select()
